# Update hashcode values (column B) in the hashcode.csv sheet
# as part of the automatic hashcode metadata refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @(
    @{Row = 43; Value = "8731886ad1651b7e01e307ed4a8e7b13"},
    @{Row = 59; Value = "007885684afcd87e61b835b2eff5ec66"},
    @{Row = 62; Value = "bf3f825dace7f973f5c47f63a384c3aa"},
    @{Row = 84; Value = "89275cf0f938689bf0b1fab22932a1ec"},
    @{Row = 86; Value = "4f28fd090ec70009da7912517f4d7939"},
    @{Row = 88; Value = "6a2a4cb969fc018679d8d02ef91fa4fd"},
    @{Row = 109; Value = "2feb3e430d71b3e18af11f8c06475e73"},
    @{Row = 115; Value = "5988191e18968528c0e3508b40458aa5"},
    @{Row = 134; Value = "7862797fb418894a33d0f6b8c9e3d362"},
    @{Row = 148; Value = "39692ec425121555be5a28fe9843bfed"},
    @{Row = 186; Value = "37917aa72ce6b7e6787d3122ec526360"},
    @{Row = 191; Value = "295819ab10107e5b676516d3e1b806e6"},
    @{Row = 195; Value = "4d5c83d5e376710af3a7a0b638681012"},
    @{Row = 196; Value = "92be4f66dcb35b1aabb9e9ec15b59464"},
    @{Row = 214; Value = "5dc85b75f4f052d7c9ad2f65403138c6"},
    @{Row = 318; Value = "3a6fc6c31b69d214d208b3c8c45b27af"},
    @{Row = 338; Value = "f4780a321a5d9da2b3be36ddb3bd1984"},
    @{Row = 340; Value = "68c23e4884abfb1affe14b6b692fc252"},
    @{Row = 348; Value = "6072a4f7a5943c4940f05be9fe6c16d0"},
    @{Row = 357; Value = "222dbf71ac6c66c8c714e52aace06047"},
    @{Row = 359; Value = "876287e2ba04c88efa3983ce3193e4c5"},
    @{Row = 363; Value = "79f9bff9e35bccedf3ba3db7bed46ac2"},
    @{Row = 380; Value = "545ba14e78d9a71615aeef7dc9dd072b"},
    @{Row = 390; Value = "ac97649e59ade6a70b20cce92ced7277"},
    @{Row = 392; Value = "1bd1ccc68993fef1d9005008d00337a6"},
    @{Row = 395; Value = "da5f629e01a0dbd5a65d4b287a3f244e"},
    @{Row = 411; Value = "e90a87636e7887f2de36a5ea0a85fe88"},
    @{Row = 431; Value = "100e4a78d4e450c1fc459d055e746e42"},
    @{Row = 439; Value = "9a2a0275ed3a403d850610246a634cd1"},
    @{Row = 459; Value = "6262ac7b545b0e78f839327ae772e388"},
    @{Row = 480; Value = "0fe57b3149dac462344231936f3e459f"},
    @{Row = 481; Value = "0543cccb29e02658d8be5593f832b5c0"},
    @{Row = 485; Value = "56b333c299c223ddd48662003bb6079d"},
    @{Row = 505; Value = "b8c6a6ec2cdafb370334752e2e75343d"},
    @{Row = 523; Value = "a07105a1fa711bd4d05e06e5467f7e49"},
    @{Row = 537; Value = "072683212b81c0658374ca875bfc979c"},
    @{Row = 555; Value = "b090b0ce8bd23ffc8f5d35d7f812ad27"},
    @{Row = 599; Value = "be3ba1b919df7e6ec6a093b13c0ac6f6"},
    @{Row = 652; Value = "d7b6be59956c29a328122c20c93bb606"},
    @{Row = 653; Value = "a7ad7fe493a9c90c26cf15b449ec2280"},
    @{Row = 683; Value = "6bb33bf273e51b34fa06e8adec039dc5"},
    @{Row = 706; Value = "3f3ea9f55fc3c56f8ef6a46339978478"},
    @{Row = 709; Value = "eaa63d292a339a16be917d4bf8677b30"},
    @{Row = 717; Value = "6678c5dfa1bc024bb15374f1b15324c9"},
    @{Row = 719; Value = "dd89cbff9a3eee2943246d621611071d"},
    @{Row = 732; Value = "da140a8fa16fbe07a90c2103e0fe6742"},
    @{Row = 745; Value = "617821a612cd696275c262f9bf9a5c54"},
    @{Row = 747; Value = "fa84e21f05f535c6ed6384d2df3fed9b"},
    @{Row = 764; Value = "bbe7ae1ac73de11278739e61e22a38c5"},
    @{Row = 801; Value = "fe192ab96888a12146d0672b606e31d5"},
    @{Row = 820; Value = "8edfce077b6ef372fd829ff5f37a55be"},
    @{Row = 823; Value = "689585849bf42fb7b5253bb978022c08"},
    @{Row = 899; Value = "da6b92ff603ee7cf879c8670684cd946"},
    @{Row = 912; Value = "6f4532b8455f950f2eab8425f40e66f5"},
    @{Row = 921; Value = "8d192c6c6c05b2945fee903f4aeb2db2"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.Value
}
